$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 91.666664
$ws.Range("I4").Value = 91.666664
$ws.Range("K4").Value = 91.666664
$ws.Range("M4").Value = 22.333336
$ws.Range("H28").Value = 467.53333
$ws.Range("I28").Value = 359.41666
$ws.Range("K28").Value = 359.41666
$ws.Range("M28").Value = 125.58334
$ws.Range("H34").Value = 12057.333
$ws.Range("I34").Value = 861
$ws.Range("J34").Value = 34450
$ws.Range("K34").Value = 861
$ws.Range("L34").Value = 34450
$ws.Range("M34").Value = -658
$ws.Range("N34").Value = -34856
$ws.Range("H36").Value = 12057.333
$ws.Range("I36").Value = 861
$ws.Range("J36").Value = 34450
$ws.Range("K36").Value = 861
$ws.Range("L36").Value = 34450
$ws.Range("M36").Value = -146
$ws.Range("N36").Value = -35880
$ws.Range("H53").Value = 136.42105
$ws.Range("I53").Value = 80.78570999999999
$ws.Range("J53").Value = 292.2
$ws.Range("K53").Value = 80.78570999999999
$ws.Range("L53").Value = 292.2
$ws.Range("M53").Value = 556.21429
$ws.Range("N53").Value = -1566.2
$ws.Range("H62").Value = 11729
$ws.Range("I62").Value = 1898.75
$ws.Range("J62").Value = 51050
$ws.Range("K62").Value = 1898.75
$ws.Range("L62").Value = 51050
$ws.Range("M62").Value = -1274.75
$ws.Range("N62").Value = -52298
$ws.Range("H64").Value = 3008.9119
$ws.Range("J64").Value = 3186.8667
$ws.Range("L64").Value = 3186.8667
$ws.Range("N64").Value = -3682.8667
$ws.Range("H65").Value = 11729
$ws.Range("I65").Value = 1898.75
$ws.Range("J65").Value = 51050
$ws.Range("K65").Value = 9493.75
$ws.Range("L65").Value = 255250
$ws.Range("M65").Value = -6373.75
$ws.Range("N65").Value = -261490
$ws.Range("H67").Value = 3008.9119
$ws.Range("J67").Value = 3186.8667
$ws.Range("L67").Value = 3186.8667
$ws.Range("N67").Value = -4902.8667
$ws.Range("H92").Value = 924.875
$ws.Range("I92").Value = 917.3333
$ws.Range("J92").Value = 947.5
$ws.Range("K92").Value = 917.3333
$ws.Range("L92").Value = 947.5
$ws.Range("M92").Value = 330.6667
$ws.Range("N92").Value = -3443.5
$ws.Range("H98").Value = 3475.75
$ws.Range("I98").Value = 2000
$ws.Range("J98").Value = 3686.5715
$ws.Range("K98").Value = 2000
$ws.Range("L98").Value = 3686.5715
$ws.Range("M98").Value = -502
$ws.Range("N98").Value = -6682.5715
$ws.Range("H107").Value = 200562
$ws.Range("I107").Value = 250452.5
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 250452.5
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -248532.5
$ws.Range("N107").Value = -4840
$ws.Range("H115").Value = 595
$ws.Range("I115").Value = 293
$ws.Range("J115").Value = 1350
$ws.Range("K115").Value = 879
$ws.Range("L115").Value = 4050
$ws.Range("M115").Value = 688
$ws.Range("N115").Value = -7184
$ws.Range("H116").Value = 3657.625
$ws.Range("I116").Value = 1725.7916
$ws.Range("J116").Value = 6555.375
$ws.Range("K116").Value = 1725.7916
$ws.Range("L116").Value = 6555.375
$ws.Range("M116").Value = 1716.2084
$ws.Range("N116").Value = -13439.375
$ws.Range("H122").Value = 3475.75
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 3686.5715
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 11059.7145
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -15959.7145
$ws.Range("H132").Value = 2653.0625
$ws.Range("I132").Value = 1406.075
$ws.Range("J132").Value = 8888
$ws.Range("K132").Value = 4218.225
$ws.Range("L132").Value = 26664
$ws.Range("M132").Value = -1688.225
$ws.Range("N132").Value = -31724
$ws.Range("H138").Value = 2374.0134
$ws.Range("I138").Value = 2104.353
$ws.Range("J138").Value = 2453.0518
$ws.Range("K138").Value = 6313.059
$ws.Range("L138").Value = 7359.155400000001
$ws.Range("M138").Value = -1173.059
$ws.Range("N138").Value = -17639.1554

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1560.6945
$ws.Range("I122").Value = 1424.871
$ws.Range("J122").Value = 2402.8
$ws.Range("K122").Value = 4274.613
$ws.Range("L122").Value = 7208.400000000001
$ws.Range("M122").Value = -1824.613
$ws.Range("N122").Value = -12108.4
$ws.Range("H123").Value = 664940
$ws.Range("J123").Value = 664940
$ws.Range("L123").Value = 664940
$ws.Range("N123").Value = -674740

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 20060
$ws.Range("J60").Value = 20060
$ws.Range("L60").Value = 20060
$ws.Range("N60").Value = -21258

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3939.1538
$ws.Range("I31").Value = 2090.2222
$ws.Range("J31").Value = 8099.25
$ws.Range("K31").Value = 2090.2222
$ws.Range("L31").Value = 8099.25
$ws.Range("M31").Value = -1795.2222
$ws.Range("N31").Value = -8689.25
$ws.Range("H34").Value = 3939.1538
$ws.Range("I34").Value = 2090.2222
$ws.Range("J34").Value = 8099.25
$ws.Range("K34").Value = 2090.2222
$ws.Range("L34").Value = 8099.25
$ws.Range("M34").Value = -1888.2222
$ws.Range("N34").Value = -8503.25
$ws.Range("H53").Value = 44499.75
$ws.Range("J53").Value = 44499.75
$ws.Range("L53").Value = 44499.75
$ws.Range("N53").Value = -45713.75
$ws.Range("H94").Value = 1150.6428
$ws.Range("I94").Value = 1044
$ws.Range("J94").Value = 1342.6
$ws.Range("K94").Value = 1044
$ws.Range("L94").Value = 1342.6
$ws.Range("M94").Value = -593
$ws.Range("N94").Value = -2244.6
$ws.Range("H107").Value = 1188.1904
$ws.Range("I107").Value = 1055
$ws.Range("J107").Value = 1309.2727
$ws.Range("K107").Value = 1055
$ws.Range("L107").Value = 1309.2727
$ws.Range("M107").Value = 865
$ws.Range("N107").Value = -5149.2727
$ws.Range("H122").Value = 2939.5908
$ws.Range("I122").Value = 1707.1111
$ws.Range("K122").Value = 5121.3333
$ws.Range("M122").Value = -2671.3333

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 491.6
$ws.Range("I5").Value = 455.2857
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1365.8571
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1253.8571
$ws.Range("N5").Value = -3224
$ws.Range("H122").Value = 3736.848
$ws.Range("I122").Value = 433.33334
$ws.Range("J122").Value = 3967.3257
$ws.Range("K122").Value = 3900.00006
$ws.Range("L122").Value = 35705.9313
$ws.Range("M122").Value = -1450.00006
$ws.Range("N122").Value = -40605.9313
$ws.Range("H135").Value = 491.6
$ws.Range("I135").Value = 455.2857
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 4097.571300000001
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -1562.571300000001
$ws.Range("N135").Value = -14070

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 673.37036
$ws.Range("I107").Value = 468
$ws.Range("J107").Value = 864.0714
$ws.Range("K107").Value = 468
$ws.Range("L107").Value = 864.0714
$ws.Range("M107").Value = 1452
$ws.Range("N107").Value = -4704.0714
$ws.Range("H122").Value = 3797.6
$ws.Range("I122").Value = 3742.182
$ws.Range("J122").Value = 3950
$ws.Range("K122").Value = 11226.546
$ws.Range("L122").Value = 11850
$ws.Range("M122").Value = -8776.545999999998
$ws.Range("N122").Value = -16750
$ws.Range("H123").Value = 16969.777
$ws.Range("J123").Value = 16969.777
$ws.Range("L123").Value = 16969.777
$ws.Range("N123").Value = -21869.777

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1356.375
$ws.Range("I46").Value = 1218.3636
$ws.Range("J46").Value = 1660
$ws.Range("K46").Value = 1218.3636
$ws.Range("L46").Value = 1660
$ws.Range("M46").Value = -1030.3636
$ws.Range("N46").Value = -2036
$ws.Range("H100").Value = 1850
$ws.Range("I100").Value = 1850
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1850
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1309
$ws.Range("N100").Value = $null

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1443.8948
$ws.Range("I107").Value = 1010.3077
$ws.Range("J107").Value = 2383.3333
$ws.Range("K107").Value = 3030.9231
$ws.Range("L107").Value = 7149.999899999999
$ws.Range("M107").Value = -1110.9231
$ws.Range("N107").Value = -10989.9999
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H113").Value = 39451.348
$ws.Range("I113").Value = 71641.57000000001
$ws.Range("J113").Value = 1896.0834
$ws.Range("K113").Value = 214924.71
$ws.Range("L113").Value = 5688.2502
$ws.Range("M113").Value = -212754.71
$ws.Range("N113").Value = -10028.2502
$ws.Range("H133").Value = 87980
$ws.Range("J133").Value = 87980
$ws.Range("L133").Value = 87980
$ws.Range("N133").Value = -98100
